$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1923076923076923
$ws.Range("C2").Value = 0.5555555555555556
$ws.Range("J2").Value = 0.008547008547008548
$ws.Range("P2").Value = 0.141025641025641
$ws.Range("S2").Value = 0.1025641025641026
$ws.Range("B3").Value = 0.007246376811594203
$ws.Range("C3").Value = 0.05072463768115942
$ws.Range("J3").Value = 0.02898550724637681
$ws.Range("P3").Value = 0.7898550724637681
$ws.Range("S3").Value = 0.1231884057971015
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6470588235294118
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.08808290155440414
$ws.Range("F6").Value = 0.09844559585492228
$ws.Range("J6").Value = 0.1606217616580311
$ws.Range("O6").Value = 0.02590673575129534
$ws.Range("Q6").Value = 0.1398963730569948
$ws.Range("R6").Value = 0.06735751295336788
$ws.Range("S6").Value = 0.4196891191709844
$ws.Range("B7").Value = 0.08433734939759036
$ws.Range("D7").Value = 0.02409638554216868
$ws.Range("F7").Value = 0.03012048192771084
$ws.Range("J7").Value = 0.09036144578313253
$ws.Range("O7").Value = 0.01204819277108434
$ws.Range("Q7").Value = 0.2048192771084337
$ws.Range("R7").Value = 0.06626506024096386
$ws.Range("S7").Value = 0.4879518072289157
$ws.Range("B8").Value = 0.08997955010224949
$ws.Range("D8").Value = 0.02249488752556237
$ws.Range("E8").Value = 0.00408997955010225
$ws.Range("F8").Value = 0.06748466257668712
$ws.Range("J8").Value = 0.1022494887525562
$ws.Range("O8").Value = 0.01431492842535787
$ws.Range("Q8").Value = 0.1860940695296524
$ws.Range("R8").Value = 0.07566462167689161
$ws.Range("S8").Value = 0.4376278118609407
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.02116402116402116
$ws.Range("F9").Value = 0.06349206349206349
$ws.Range("J9").Value = 0.1005291005291005
$ws.Range("O9").Value = 0.02645502645502645
$ws.Range("Q9").Value = 0.1587301587301587
$ws.Range("R9").Value = 0.06878306878306878
$ws.Range("S9").Value = 0.4656084656084656
$ws.Range("B10").Value = 0.09225092250922509
$ws.Range("D10").Value = 0.01383763837638376
$ws.Range("F10").Value = 0.06549815498154982
$ws.Range("J10").Value = 0.1060885608856089
$ws.Range("O10").Value = 0.01383763837638376
$ws.Range("Q10").Value = 0.2177121771217712
$ws.Range("R10").Value = 0.07564575645756458
$ws.Range("S10").Value = 0.4151291512915129
$ws.Range("G11").Value = 0.1654676258992806
$ws.Range("J11").Value = 0.06474820143884892
$ws.Range("K11").Value = 0.2122302158273381
$ws.Range("L11").Value = 0.5215827338129496
$ws.Range("S11").Value = 0.03597122302158273
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.1790123456790123
$ws.Range("K12").Value = 0.02469135802469136
$ws.Range("L12").Value = 0.04320987654320987
$ws.Range("S12").Value = 0.08641975308641975
$ws.Range("G13").Value = 0.5945945945945946
$ws.Range("J13").Value = 0.2432432432432433
$ws.Range("S13").Value = 0.1621621621621622
$ws.Range("F15").Value = 0.02242152466367713
$ws.Range("H15").Value = 0.1838565022421525
$ws.Range("I15").Value = 0.05829596412556054
$ws.Range("J15").Value = 0.3587443946188341
$ws.Range("K15").Value = 0.07174887892376682
$ws.Range("M15").Value = 0.008968609865470852
$ws.Range("O15").Value = 0.09865470852017937
$ws.Range("S15").Value = 0.1973094170403587
$ws.Range("H16").Value = 0.286624203821656
$ws.Range("I16").Value = 0.08917197452229299
$ws.Range("J16").Value = 0.3503184713375796
$ws.Range("K16").Value = 0.09554140127388536
$ws.Range("M16").Value = 0.01910828025477707
$ws.Range("N16").Value = 0.006369426751592357
$ws.Range("O16").Value = 0.05732484076433121
$ws.Range("S16").Value = 0.09554140127388536
$ws.Range("F17").Value = 0.009685230024213076
$ws.Range("H17").Value = 0.1719128329297821
$ws.Range("I17").Value = 0.1186440677966102
$ws.Range("J17").Value = 0.3898305084745763
$ws.Range("K17").Value = 0.07990314769975787
$ws.Range("M17").Value = 0.01210653753026634
$ws.Range("O17").Value = 0.05569007263922518
$ws.Range("S17").Value = 0.162227602905569
$ws.Range("F18").Value = 0.01290322580645161
$ws.Range("H18").Value = 0.1419354838709677
$ws.Range("I18").Value = 0.1096774193548387
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.07096774193548387
$ws.Range("M18").Value = 0.01935483870967742
$ws.Range("O18").Value = 0.05161290322580645
$ws.Range("S18").Value = 0.1935483870967742
$ws.Range("F19").Value = 0.01641791044776119
$ws.Range("H19").Value = 0.2335820895522388
$ws.Range("I19").Value = 0.0664179104477612
$ws.Range("J19").Value = 0.3365671641791045
$ws.Range("K19").Value = 0.1
$ws.Range("M19").Value = 0.01716417910447761
$ws.Range("O19").Value = 0.07313432835820896
$ws.Range("S19").Value = 0.1567164179104478
